$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Version bump: 3.0.0 -> 3.7.0 everywhere it appears as a literal
#    version string (git examples, Tomcat <Context> docBase lines,
#    the war-file backup/move commands).
# ------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("3.0.0", $false, $false, $false, $false, $false, $true, 1, $false, "3.7.0", 2)
Write-Output "version bump found=$found1"

# ------------------------------------------------------------------
# 2) Tidy the "stop" sentence: re-assert the same text over the
#    existing range so the run split collapses into a single run
#    (this also drops the now-redundant grammar-check markers around
#    the word "stop").
# ------------------------------------------------------------------
$r2 = $d.Content
$stopText = "If the web application URL displays Tomcat home page, stop and start Apache Tomcat "
$found2 = $r2.Find.Execute($stopText, $false, $false, $false, $false, $false, $true, 1, $false, $stopText, 2)
Write-Output "stop sentence found=$found2"

# ------------------------------------------------------------------
# 3) Add the Java 21 module-system flag to the setenv.sh JAVA_OPTS
#    snippet, right after the logback configurationFile option and
#    before the closing quote.
# ------------------------------------------------------------------
$r3 = $d.Content
$anchor = 'Dlogback.configurationFile=$CATALINA_HOME/conf/hpc-server/logback.xml'
$found3 = $r3.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "logback anchor found=$found3"
if ($found3) {
    $insertPoint = $r3.Duplicate
    $insertPoint.Collapse(0)
    $insertPoint.InsertAfter(" --add-opens java.base/java.net=ALL-UNNAMED")
}
